$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the title
# (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------

$titlePar = $d.Paragraphs.Item(1)
$titlePar.Range.InsertParagraphAfter()

# Copy the formatted contents (incl. the leading empty run artifact) of the
# "The graphics of ..." paragraph as a scaffold for our new paragraph, then
# overwrite its text piece by piece so the leading empty run survives.
$scaffoldPar = $d.Paragraphs.Item(4)
$scaffoldFT = $scaffoldPar.Range.FormattedText
$d.Paragraphs.Item(2).Range.FormattedText = $scaffoldFT

$metaPar = $d.Paragraphs.Item(2)
$metaStart = $metaPar.Range.Start
$metaFullText = $metaPar.Range.Text
$italicIdx = $metaFullText.IndexOf("9 Burning Stars")
$tailIdx = $metaFullText.IndexOf(" are out of this world")

$metaLabel = "Meta description"

# Bold the first run *before* touching its text so it doesn't get merged
# into the empty run that precedes it.
$labelRange = $d.Range($metaStart, $metaStart + $italicIdx)
$labelRange.Font.Bold = 1
$labelRange.Text = $metaLabel

# Drop the (now stale) italicised "9 Burning Stars" run text - this merges
# it with the plain trailing run that follows it.
$afterLabelStart = $metaStart + $metaLabel.Length
$refreshedText = $d.Paragraphs.Item(2).Range.Text
$tailIdx2 = $refreshedText.IndexOf(" are out of this world")
$italicRange = $d.Range($afterLabelStart, $metaStart + $tailIdx2)
$italicRange.Text = ""

# Replace the remaining (plain formatted) run with the real meta description
# text.
$descPar = $d.Paragraphs.Item(2)
$descStart = $metaStart + $metaLabel.Length
$descEnd = $descPar.Range.End - 1
$descRange = $d.Range($descStart, $descEnd)
$descRange.Text = ": Read a review of the slot game 9 Burning Stars, including pros and cons. Try it for free and enjoy the elaborate graphics and engaging music."

# The paragraph inherited the Heading1 style from its predecessor; put it
# back to the document's normal body-text style.
$d.Paragraphs.Item(2).Style = "Normal"

# ---------------------------------------------------------------------------
# Change 2: drop the duplicated bold "Play 9 Burning Stars..." paragraph near
# the end of the document, and replace the italic paragraph's text with the
# feature-image prompt.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupTitlePar = $d.Paragraphs.Item($count - 1)
$dupTitlePar.Range.Delete()

$count2 = $d.Paragraphs.Count
$imgPar = $d.Paragraphs.Item($count2)
$imgStart = $imgPar.Range.Start
$imgEnd = $imgPar.Range.End - 1
$imgRange = $d.Range($imgStart, $imgEnd)
$imgRange.Text = "For the feature image fitting the game `"9 Burning Stars`", please create a cartoon-style image that features a happy Maya warrior with glasses. The warrior should be holding a torch with fire emanating from it, while standing in front of a starry galaxy background. The warrior should be depicted as excited and triumphant with vibrant energy surrounding them, symbolizing the thrill of playing the game. The image should include game elements such as the logo, the chili symbol representing volatility choice, and the burning star symbol. The image should also include a tagline such as `"Experience the Galactic Thrill with 9 Burning Stars`". The overall image should convey the excitement, energy, and fun of playing the game, while highlighting its unique theme and features."

Write-Output "done"
